$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the
#    H1 title paragraph ("Play Deep Sea Magic Slot Game for Free - Review").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaText = "Meta description: Find out about the features of the Deep Sea Magic slot game and play it for free. Read our review before playing for real money."
$metaPara.Range.Text = $metaText

# Bold just the "Meta description" label (the first 17 characters).
$labelStart = $metaPara.Range.Start
$labelEnd = $labelStart + "Meta description".Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Near the end of the document, the title/description pair was
#    duplicated. Remove the second (bold) "Play Deep Sea Magic Slot
#    Game for Free - Review" paragraph entirely, keeping only the
#    italic paragraph, whose text becomes a DALLE image prompt.
# ------------------------------------------------------------------
$searchRng = $d.Content
$searchRng.Find.Text = "Play Deep Sea Magic Slot Game for Free - Review"
$searchRng.Find.Forward = $true
$searchRng.Find.Wrap = 0

$null = $searchRng.Find.Execute()          # 1st hit: the H1 title itself
$found = $searchRng.Find.Execute()         # 2nd hit: the duplicated paragraph

if ($found -and $searchRng.Start -gt 0) {
    $dupPara = $searchRng.Paragraphs(1)
    $dupPara.Range.Delete()
}

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$oldDesc = "Find out about the features of the Deep Sea Magic slot game and play it for free. Read our review before playing for real money."
$newDesc = "DALLE, please create a feature image for `"Deep Sea Magic`" that fits the game's theme and features a happy Maya warrior with glasses in a cartoon style. The image should be eye-catching and playful, incorporating elements of the deep sea and the game's features such as the Drop & Lock feature, bonuses, and wild symbols. Please make sure that the image is high-quality and in line with the overall aesthetic of the game. Thank you!"

# The very last paragraph in the document behaves specially: assigning
# to its own cached .Range.Text inserts rather than replaces. Re-build
# the range from the Document so the replace actually overwrites the
# old text instead of prepending the new text in front of it.
$descStart = $lastPara.Range.Start
$descEnd = $lastPara.Range.End
if ($lastPara.Range.Text.EndsWith([char]13)) {
    $descEnd = $descEnd - 1
}
$descRange = $d.Range($descStart, $descEnd)

if ($descRange.Text -eq $oldDesc) {
    $descRange.Text = $newDesc
}
